$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hidden_5")
$cols = $ws.Range("B1:W1").EntireColumn
$cols.ColumnWidth = 9.140625
